$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 38245.03311329238
$ws.Range("C2").Value = 1757.557135943058
$ws.Range("D2").Value = 6200.304207411717
$ws.Range("E2").Value = 2381.49433727492
$ws.Range("F2").Value = 221.5532468662895
$ws.Range("G2").Value = 847.4837575814532
$ws.Range("H2").Value = 40.89049071784314
$ws.Range("I2").Value = 3228.978094856374
$ws.Range("J2").Value = 262.4437375841327
$ws.Range("K2").Value = 384.5313542794729
$ws.Range("L2").Value = 25.17144328777763
$ws.Range("M2").Value = 276.9310007360748
$ws.Range("N2").Value = 39.70680705549299
$ws.Range("O2").Value = 31440.17263207441
$ws.Range("P2").Value = 1891.126053113459
$ws.Range("Q2").Value = 6189.874009047865
$ws.Range("R2").Value = 2532.094338317329
$ws.Range("S2").Value = 222.136886035803
$ws.Range("T2").Value = 621.4769208533244
$ws.Range("U2").Value = 31.67411357577021
$ws.Range("V2").Value = 3153.571259170653
$ws.Range("W2").Value = 253.8109996115732
$ws.Range("X2").Value = 343.6188554092861
$ws.Range("Y2").Value = 19.35495179734423
$ws.Range("Z2").Value = 246.6065971555531
$ws.Range("AA2").Value = 43.64389993789622

$ws.Range("B3").Value = 32499.04413726227
$ws.Range("C3").Value = 1786.609670827545
$ws.Range("D3").Value = 5378.734706132076
$ws.Range("E3").Value = 2007.90645452699
$ws.Range("F3").Value = 195.3979334832982
$ws.Range("G3").Value = 617.2524073440741
$ws.Range("H3").Value = 33.34520400626526
$ws.Range("I3").Value = 2625.158861871064
$ws.Range("J3").Value = 228.7431374895635
$ws.Range("K3").Value = 336.8281547835174
$ws.Range("L3").Value = 21.20098980071841
$ws.Range("M3").Value = 234.4034127069872
$ws.Range("N3").Value = 37.37377923268704
$ws.Range("O3").Value = 36289.31985932789
$ws.Range("P3").Value = 1900.827443428021
$ws.Range("Q3").Value = 7294.948835296573
$ws.Range("R3").Value = 2551.357721439193
$ws.Range("S3").Value = 234.3660304445436
$ws.Range("T3").Value = 566.7822713398243
$ws.Range("U3").Value = 31.4896106465157
$ws.Range("V3").Value = 3118.139992779018
$ws.Range("W3").Value = 265.8556410910593
$ws.Range("X3").Value = 382.8566004433583
$ws.Range("Y3").Value = 19.54214375698277
$ws.Range("Z3").Value = 281.3605260728197
$ws.Range("AA3").Value = 47.78193904311267

$ws.Range("B4").Value = 39835.65461768988
$ws.Range("C4").Value = 1642.562416926959
$ws.Range("D4").Value = 5888.771798129625
$ws.Range("E4").Value = 2214.43062076898
$ws.Range("F4").Value = 206.5063738952491
$ws.Range("G4").Value = 651.4409919283753
$ws.Range("H4").Value = 34.96197005383755
$ws.Range("I4").Value = 2865.871612697355
$ws.Range("J4").Value = 241.4683439490866
$ws.Range("K4").Value = 384.7696365088916
$ws.Range("L4").Value = 20.71681271646442
$ws.Range("M4").Value = 288.9744933766121
$ws.Range("N4").Value = 39.26688892702936
$ws.Range("O4").Value = 33059.42705334981
$ws.Range("P4").Value = 1710.777219034276
$ws.Range("Q4").Value = 5854.851126469145
$ws.Range("R4").Value = 2330.013959716519
$ws.Range("S4").Value = 213.6248492600263
$ws.Range("T4").Value = 619.7484351617462
$ws.Range("U4").Value = 30.44134379420773
$ws.Range("V4").Value = 2949.762394878265
$ws.Range("W4").Value = 244.0661930542341
$ws.Range("X4").Value = 366.0246009398963
$ws.Range("Y4").Value = 17.43523994079843
$ws.Range("Z4").Value = 268.7016179673131
$ws.Range("AA4").Value = 42.72054323117423

$ws.Range("B5").Value = 36329.3198039069
$ws.Range("C5").Value = 1796.3453528688
$ws.Range("D5").Value = 5485.629574808951
$ws.Range("E5").Value = 2253.566197583306
$ws.Range("F5").Value = 214.0807679957353
$ws.Range("G5").Value = 831.1872437403199
$ws.Range("H5").Value = 41.92023953273372
$ws.Range("I5").Value = 3084.753441323626
$ws.Range("J5").Value = 256.001007528469
$ws.Range("K5").Value = 361.8716943300348
$ws.Range("L5").Value = 22.9985057976689
$ws.Range("M5").Value = 260.5669865292919
$ws.Range("N5").Value = 39.72274921270161
$ws.Range("O5").Value = 35826.4890328753
$ws.Range("P5").Value = 1873.92649046366
$ws.Range("Q5").Value = 6339.066088984951
$ws.Range("R5").Value = 2402.213546170371
$ws.Range("S5").Value = 231.0500634971331
$ws.Range("T5").Value = 679.3546493931528
$ws.Range("U5").Value = 34.26146077844262
$ws.Range("V5").Value = 3081.568195563524
$ws.Range("W5").Value = 265.3115242755757
$ws.Range("X5").Value = 404.7662379128714
$ws.Range("Y5").Value = 21.82375188926511
$ws.Range("Z5").Value = 291.1439333807644
$ws.Range("AA5").Value = 52.40385255022856

$ws.Range("B6").Value = 38364.30976201847
$ws.Range("C6").Value = 1813.472520353595
$ws.Range("D6").Value = 6443.928726638457
$ws.Range("E6").Value = 2354.562312141525
$ws.Range("F6").Value = 223.5256608571924
$ws.Range("G6").Value = 841.4474052927986
$ws.Range("H6").Value = 42.20616749958173
$ws.Range("I6").Value = 3196.009717434323
$ws.Range("J6").Value = 265.7318283567741
$ws.Range("K6").Value = 387.1249594865291
$ws.Range("L6").Value = 23.39365086153757
$ws.Range("M6").Value = 276.3745637222472
$ws.Range("N6").Value = 41.63703451013856
$ws.Range("O6").Value = 32315.64799448937
$ws.Range("P6").Value = 1895.49159348734
$ws.Range("Q6").Value = 6329.355380213569
$ws.Range("R6").Value = 2522.849394956393
$ws.Range("S6").Value = 232.5703907408471
$ws.Range("T6").Value = 683.3895994405326
$ws.Range("U6").Value = 34.39571636840808
$ws.Range("V6").Value = 3206.238994396926
$ws.Range("W6").Value = 266.9661071092552
$ws.Range("X6").Value = 365.9224841371235
$ws.Range("Y6").Value = 19.42781126609105
$ws.Range("Z6").Value = 257.7601949884595
$ws.Range("AA6").Value = 47.46906644811116

$ws.Range("B7").Value = 35727.34822351123
$ws.Range("C7").Value = 1824.941744244175
$ws.Range("D7").Value = 5778.054733501265
$ws.Range("E7").Value = 2374.034468055822
$ws.Range("F7").Value = 224.3303883474233
$ws.Range("G7").Value = 862.4096956328447
$ws.Range("H7").Value = 44.17597803273693
$ws.Range("I7").Value = 3236.444163688667
$ws.Range("J7").Value = 268.5063663801603
$ws.Range("K7").Value = 373.9988628839857
$ws.Range("L7").Value = 22.42705974676062
$ws.Range("M7").Value = 266.3084704183372
$ws.Range("N7").Value = 39.84229381419119
$ws.Range("O7").Value = 35631.23939376771
$ws.Range("P7").Value = 1889.975864980121
$ws.Range("Q7").Value = 6538.871685802329
$ws.Range("R7").Value = 2553.280571057361
$ws.Range("S7").Value = 244.9416110812896
$ws.Range("T7").Value = 685.0519844463673
$ws.Range("U7").Value = 37.24662424277886
$ws.Range("V7").Value = 3238.332555503728
$ws.Range("W7").Value = 282.1882353240684
$ws.Range("X7").Value = 404.2059888197002
$ws.Range("Y7").Value = 20.576745637376
$ws.Range("Z7").Value = 290.8971981268087
$ws.Range("AA7").Value = 49.26548383392851

$ws.Range("B8").Value = 34951.32028616372
$ws.Range("C8").Value = 1751.151805984811
$ws.Range("D8").Value = 5539.375957080137
$ws.Range("E8").Value = 2214.672902506505
$ws.Range("F8").Value = 211.3902208503606
$ws.Range("G8").Value = 840.4830441363794
$ws.Range("H8").Value = 42.23366861559504
$ws.Range("I8").Value = 3055.155946642884
$ws.Range("J8").Value = 253.6238894659556
$ws.Range("K8").Value = 363.8326167412018
$ws.Range("L8").Value = 21.65071469694586
$ws.Range("M8").Value = 255.69757874504
$ws.Range("N8").Value = 37.9413581141843
$ws.Range("O8").Value = 36022.98809504891
$ws.Range("P8").Value = 1771.447233404109
$ws.Range("Q8").Value = 6460.413959497083
$ws.Range("R8").Value = 2495.709445324508
$ws.Range("S8").Value = 238.0490047412358
$ws.Range("T8").Value = 658.3429792237313
$ws.Range("U8").Value = 35.26213437837964
$ws.Range("V8").Value = 3154.052424548239
$ws.Range("W8").Value = 273.3111391196155
$ws.Range("X8").Value = 413.4399324186161
$ws.Range("Y8").Value = 20.2577579778538
$ws.Range("Z8").Value = 300.2766371495229
$ws.Range("AA8").Value = 51.06827590680861

$ws.Range("B9").Value = 31712.65601471566
$ws.Range("C9").Value = 1697.513342961181
$ws.Range("D9").Value = 4728.006377961241
$ws.Range("E9").Value = 1952.333528717393
$ws.Range("F9").Value = 192.5132950778745
$ws.Range("G9").Value = 704.3459205458515
$ws.Range("H9").Value = 36.80827451024174
$ws.Range("I9").Value = 2656.679449263244
$ws.Range("J9").Value = 229.3215695881162
$ws.Range("K9").Value = 331.0493211637479
$ws.Range("L9").Value = 18.82247802464963
$ws.Range("M9").Value = 228.5377516086306
$ws.Range("N9").Value = 38.90880828022765
$ws.Range("O9").Value = 38629.15028659082
$ws.Range("P9").Value = 1906.544449689495
$ws.Range("Q9").Value = 7057.538681514773
$ws.Range("R9").Value = 2625.852751651979
$ws.Range("S9").Value = 250.3188464636979
$ws.Range("T9").Value = 700.7098024968385
$ws.Range("U9").Value = 37.74508253750755
$ws.Range("V9").Value = 3326.562554148817
$ws.Range("W9").Value = 288.0639290012055
$ws.Range("X9").Value = 433.0815423552864
$ws.Range("Y9").Value = 19.99893606800926
$ws.Range("Z9").Value = 313.2265431193254
$ws.Range("AA9").Value = 55.35588011076818

$ws.Range("B10").Value = 34599.30578078904
$ws.Range("C10").Value = 1813.768055238299
$ws.Range("D10").Value = 5636.938385052741
$ws.Range("E10").Value = 2165.221969408223
$ws.Range("F10").Value = 198.1199025707598
$ws.Range("G10").Value = 664.0143623256224
$ws.Range("H10").Value = 35.03201147873428
$ws.Range("I10").Value = 2829.236331733845
$ws.Range("J10").Value = 233.1519140494941
$ws.Range("K10").Value = 348.0469324585746
$ws.Range("L10").Value = 18.02493926423912
$ws.Range("M10").Value = 243.1907355472126
$ws.Range("N10").Value = 36.78198422600281
$ws.Range("O10").Value = 38013.70916507281
$ws.Range("P10").Value = 1934.21398782214
$ws.Range("Q10").Value = 7104.278146511607
$ws.Range("R10").Value = 2530.029499378603
$ws.Range("S10").Value = 228.0049846634972
$ws.Range("T10").Value = 607.239011985018
$ws.Range("U10").Value = 30.51344935987909
$ws.Range("V10").Value = 3137.26851136362
$ws.Range("W10").Value = 258.5184340233763
$ws.Range("X10").Value = 397.0339840138156
$ws.Range("Y10").Value = 19.09787495587503
$ws.Range("Z10").Value = 294.1851427522167
$ws.Range("AA10").Value = 46.66148726448649

$ws.Range("B11").Value = 33556.87371432573
$ws.Range("C11").Value = 1814.323322454931
$ws.Range("D11").Value = 5260.150160628757
$ws.Range("E11").Value = 2196.798965504275
$ws.Range("F11").Value = 209.5475429531017
$ws.Range("G11").Value = 769.2388941444563
$ws.Range("H11").Value = 38.5938364196981
$ws.Range("I11").Value = 2966.037859648732
$ws.Range("J11").Value = 248.1413793727997
$ws.Range("K11").Value = 338.6410023512584
$ws.Range("L11").Value = 20.44505142322577
$ws.Range("M11").Value = 236.9468140317591
$ws.Range("N11").Value = 37.53058492236872
$ws.Range("O11").Value = 37795.27953249563
$ws.Range("P11").Value = 1902.70474997408
$ws.Range("Q11").Value = 6709.321416079494
$ws.Range("R11").Value = 2320.804939557562
$ws.Range("S11").Value = 221.8202145568863
$ws.Range("T11").Value = 640.4599988154947
$ws.Range("U11").Value = 31.51960780515648
$ws.Range("V11").Value = 2961.264938373056
$ws.Range("W11").Value = 253.3398223620428
$ws.Range("X11").Value = 398.2331604895535
$ws.Range("Y11").Value = 20.11182949080873
$ws.Range("Z11").Value = 300.9267137088569
$ws.Range("AA11").Value = 48.10667988009593

$ws.Range("B12").Value = 31924.77795037885
$ws.Range("C12").Value = 1709.226529569799
$ws.Range("D12").Value = 4951.300566856832
$ws.Range("E12").Value = 1989.14131314416
$ws.Range("F12").Value = 188.511910952724
$ws.Range("G12").Value = 667.0945804989935
$ws.Range("H12").Value = 34.36742642369819
$ws.Range("I12").Value = 2656.235893643154
$ws.Range("J12").Value = 222.8793373764222
$ws.Range("K12").Value = 323.1572713646102
$ws.Range("L12").Value = 19.40959774174423
$ws.Range("M12").Value = 222.1603472937751
$ws.Range("N12").Value = 33.83238760043948
$ws.Range("O12").Value = 37177.20980042957
$ws.Range("P12").Value = 1882.029628096796
$ws.Range("Q12").Value = 6997.97720980198
$ws.Range("R12").Value = 2951.406684730563
$ws.Range("S12").Value = 272.2086915090495
$ws.Range("T12").Value = 839.1136476804604
$ws.Range("U12").Value = 44.37879488398616
$ws.Range("V12").Value = 3790.520332411023
$ws.Range("W12").Value = 316.5874863930357
$ws.Range("X12").Value = 404.8806935345627
$ws.Range("Y12").Value = 21.74811024745186
$ws.Range("Z12").Value = 293.8804295773204
$ws.Range("AA12").Value = 50.92385965164715

$ws.Range("B13").Value = 33400.81349215898
$ws.Range("C13").Value = 1753.497097978196
$ws.Range("D13").Value = 5492.744281768415
$ws.Range("E13").Value = 2081.617532743563
$ws.Range("F13").Value = 191.1969619314859
$ws.Range("G13").Value = 609.0852414515706
$ws.Range("H13").Value = 32.26976470417866
$ws.Range("I13").Value = 2690.702774195134
$ws.Range("J13").Value = 223.4667266356646
$ws.Range("K13").Value = 335.7716657760849
$ws.Range("L13").Value = 18.17784094726269
$ws.Range("M13").Value = 238.7297038298212
$ws.Range("N13").Value = 37.66836162969642
$ws.Range("O13").Value = 36883.66617384075
$ws.Range("P13").Value = 1840.583584579148
$ws.Range("Q13").Value = 7100.825748611669
$ws.Range("R13").Value = 2564.479362394357
$ws.Range("S13").Value = 232.6371150823617
$ws.Range("T13").Value = 601.7674645570363
$ws.Range("U13").Value = 31.38321765866515
$ws.Range("V13").Value = 3166.246826951393
$ws.Range("W13").Value = 264.0203327410269
$ws.Range("X13").Value = 421.6225187593739
$ws.Range("Y13").Value = 18.9610097184604
$ws.Range("Z13").Value = 313.0607988831151
$ws.Range("AA13").Value = 49.10175908112218

$ws.Range("B14").Value = 32689.01741591886
$ws.Range("C14").Value = 1701.174526645647
$ws.Range("D14").Value = 5112.88799434892
$ws.Range("E14").Value = 1964.282481062862
$ws.Range("F14").Value = 187.7923593051985
$ws.Range("G14").Value = 624.7304372284598
$ws.Range("H14").Value = 32.37312598259111
$ws.Range("I14").Value = 2589.012918291322
$ws.Range("J14").Value = 220.1654852877896
$ws.Range("K14").Value = 337.2911964459875
$ws.Range("L14").Value = 18.38693421152254
$ws.Range("M14").Value = 233.117807080912
$ws.Range("N14").Value = 34.21036933884186
$ws.Range("O14").Value = 40668.16497309941
$ws.Range("P14").Value = 1853.31725983307
$ws.Range("Q14").Value = 7913.337624110826
$ws.Range("R14").Value = 2996.468541485473
$ws.Range("S14").Value = 282.3192247311213
$ws.Range("T14").Value = 825.5512108868013
$ws.Range("U14").Value = 43.97405884210976
$ws.Range("V14").Value = 3822.019752372274
$ws.Range("W14").Value = 326.2932835732311
$ws.Range("X14").Value = 443.1763366070703
$ws.Range("Y14").Value = 21.16522445671716
$ws.Range("Z14").Value = 325.5251420143118
$ws.Range("AA14").Value = 52.36861314931601

$ws.Range("B15").Value = 34622.12060813731
$ws.Range("C15").Value = 1798.640002538314
$ws.Range("D15").Value = 5211.75227434212
$ws.Range("E15").Value = 2058.586320347985
$ws.Range("F15").Value = 193.3570242230889
$ws.Range("G15").Value = 709.8461316887618
$ws.Range("H15").Value = 35.25963866816105
$ws.Range("I15").Value = 2768.432452036747
$ws.Range("J15").Value = 228.6166628912499
$ws.Range("K15").Value = 348.2456417825973
$ws.Range("L15").Value = 20.78726325078235
$ws.Range("M15").Value = 249.0724970614975
$ws.Range("N15").Value = 37.82702651580312
$ws.Range("O15").Value = 35956.11748919122
$ws.Range("P15").Value = 1914.32467951373
$ws.Range("Q15").Value = 6648.11881062922
$ws.Range("R15").Value = 2571.577495722455
$ws.Range("S15").Value = 241.3841052808941
$ws.Range("T15").Value = 692.7725020088776
$ws.Range("U15").Value = 36.62463872834574
$ws.Range("V15").Value = 3264.349997731333
$ws.Range("W15").Value = 278.0087440092399
$ws.Range("X15").Value = 402.5917815805749
$ws.Range("Y15").Value = 20.58818266223681
$ws.Range("Z15").Value = 295.1599233991067
$ws.Range("AA15").Value = 48.25928335638009

$ws.Range("B16").Value = 30834.49821083745
$ws.Range("C16").Value = 1824.955169983611
$ws.Range("D16").Value = 4808.138677226871
$ws.Range("E16").Value = 1966.247746387659
$ws.Range("F16").Value = 185.8928085102627
$ws.Range("G16").Value = 684.0892388843872
$ws.Range("H16").Value = 34.48313997859874
$ws.Range("I16").Value = 2650.336985272047
$ws.Range("J16").Value = 220.3759484888615
$ws.Range("K16").Value = 316.4500301182314
$ws.Range("L16").Value = 18.55154460157359
$ws.Range("M16").Value = 224.1898755806921
$ws.Range("N16").Value = 34.12680711643406
$ws.Range("O16").Value = 39333.47810827673
$ws.Range("P16").Value = 1950.234577710387
$ws.Range("Q16").Value = 7275.18438381999
$ws.Range("R16").Value = 2661.962250527356
$ws.Range("S16").Value = 242.7975147658969
$ws.Range("T16").Value = 662.3447190965047
$ws.Range("U16").Value = 32.7655537530103
$ws.Range("V16").Value = 3324.306969623861
$ws.Range("W16").Value = 275.5630685189072
$ws.Range("X16").Value = 435.2114100255834
$ws.Range("Y16").Value = 19.19037990539216
$ws.Range("Z16").Value = 327.9275164545384
$ws.Range("AA16").Value = 51.24997588561863

$ws.Range("B17").Value = 32197.91631945951
$ws.Range("C17").Value = 1805.582730617899
$ws.Range("D17").Value = 5287.461061073574
$ws.Range("E17").Value = 2182.108513671355
$ws.Range("F17").Value = 200.337664868818
$ws.Range("G17").Value = 733.223326592201
$ws.Range("H17").Value = 37.13326301539426
$ws.Range("I17").Value = 2915.331840263556
$ws.Range("J17").Value = 237.4709278842122
$ws.Range("K17").Value = 325.5494735320646
$ws.Range("L17").Value = 17.62620049982698
$ws.Range("M17").Value = 227.2553887637085
$ws.Range("N17").Value = 36.11733105354341
$ws.Range("O17").Value = 38227.21615555442
$ws.Range("P17").Value = 1907.749988971838
$ws.Range("Q17").Value = 7201.163201644224
$ws.Range("R17").Value = 2671.698051368409
$ws.Range("S17").Value = 245.2744887391809
$ws.Range("T17").Value = 703.5289181469866
$ws.Range("U17").Value = 36.12723154371011
$ws.Range("V17").Value = 3375.226969515396
$ws.Range("W17").Value = 281.4017202828911
$ws.Range("X17").Value = 411.1229662586708
$ws.Range("Y17").Value = 18.45047282502311
$ws.Range("Z17").Value = 308.5979099002224
$ws.Range("AA17").Value = 47.84697068360161

$ws.Range("B18").Value = 32464.62477656585
$ws.Range("C18").Value = 1917.084750358108
$ws.Range("D18").Value = 5720.668541323336
$ws.Range("E18").Value = 2298.195199375691
$ws.Range("F18").Value = 207.5244364210596
$ws.Range("G18").Value = 685.2584542983894
$ws.Range("H18").Value = 36.89834281183595
$ws.Range("I18").Value = 2983.453653674081
$ws.Range("J18").Value = 244.4227792328955
$ws.Range("K18").Value = 328.1650305555928
$ws.Range("L18").Value = 19.13552302802644
$ws.Range("M18").Value = 231.9727766106595
$ws.Range("N18").Value = 36.04145555465609
$ws.Range("O18").Value = 39861.06727929034
$ws.Range("P18").Value = 1990.21486077696
$ws.Range("Q18").Value = 7394.432913796735
$ws.Range("R18").Value = 2824.042567216219
$ws.Range("S18").Value = 257.8007822823
$ws.Range("T18").Value = 772.2808620996474
$ws.Range("U18").Value = 38.75414424911429
$ws.Range("V18").Value = 3596.323429315867
$ws.Range("W18").Value = 296.5549265314143
$ws.Range("X18").Value = 434.0992711116339
$ws.Range("Y18").Value = 22.32968082193736
$ws.Range("Z18").Value = 322.3668549959459
$ws.Range("AA18").Value = 53.73533695125183

$ws.Range("B19").Value = 30119.38859742292
$ws.Range("C19").Value = 1736.66476249109
$ws.Range("D19").Value = 4836.715978468995
$ws.Range("E19").Value = 1923.973250434202
$ws.Range("F19").Value = 179.8137021520933
$ws.Range("G19").Value = 611.904366424434
$ws.Range("H19").Value = 31.33086458238964
$ws.Range("I19").Value = 2535.877616858636
$ws.Range("J19").Value = 211.1445667344829
$ws.Range("K19").Value = 316.4643785870497
$ws.Range("L19").Value = 19.39716607790639
$ws.Range("M19").Value = 216.5225582828384
$ws.Range("N19").Value = 31.89491258960997
$ws.Range("O19").Value = 38094.13135671659
$ws.Range("P19").Value = 1851.720301729615
$ws.Range("Q19").Value = 6853.158811800633
$ws.Range("R19").Value = 2581.829539908847
$ws.Range("S19").Value = 252.1425839410952
$ws.Range("T19").Value = 707.9335520738163
$ws.Range("U19").Value = 37.84568357156233
$ws.Range("V19").Value = 3289.763091982663
$ws.Range("W19").Value = 289.9882675126575
$ws.Range("X19").Value = 437.8126496536024
$ws.Range("Y19").Value = 22.13134435005868
$ws.Range("Z19").Value = 321.762715523745
$ws.Range("AA19").Value = 52.6380912882653

$ws.Range("B20").Value = 30156.84747066248
$ws.Range("C20").Value = 1814.044874241087
$ws.Range("D20").Value = 4911.175893357597
$ws.Range("E20").Value = 1959.178689402509
$ws.Range("F20").Value = 183.8348164877708
$ws.Range("G20").Value = 655.7587942333535
$ws.Range("H20").Value = 33.52550380775735
$ws.Range("I20").Value = 2614.937483635863
$ws.Range("J20").Value = 217.3603202955282
$ws.Range("K20").Value = 312.2836216826745
$ws.Range("L20").Value = 19.81446849970185
$ws.Range("M20").Value = 217.6542408808185
$ws.Range("N20").Value = 35.81817550107806
$ws.Range("O20").Value = 37572.19880266338
$ws.Range("P20").Value = 1947.71611415023
$ws.Range("Q20").Value = 7220.206106041735
$ws.Range("R20").Value = 2691.761931028316
$ws.Range("S20").Value = 248.4006198059043
$ws.Range("T20").Value = 692.3382162146625
$ws.Range("U20").Value = 37.41803846101172
$ws.Range("V20").Value = 3384.100147242978
$ws.Range("W20").Value = 285.818658266916
$ws.Range("X20").Value = 414.7442705052798
$ws.Range("Y20").Value = 21.82259785782909
$ws.Range("Z20").Value = 301.9781936940815
$ws.Range("AA20").Value = 50.90108159053389

$ws.Range("B21").Value = 35156.95476618902
$ws.Range("C21").Value = 1814.791029582654
$ws.Range("D21").Value = 5968.638919531206
$ws.Range("E21").Value = 2274.696901113578
$ws.Range("F21").Value = 209.3587300746551
$ws.Range("G21").Value = 783.5486663241466
$ws.Range("H21").Value = 39.46282914811896
$ws.Range("I21").Value = 3058.245567437725
$ws.Range("J21").Value = 248.8215592227741
$ws.Range("K21").Value = 361.2164668351104
$ws.Range("L21").Value = 21.22087414459916
$ws.Range("M21").Value = 252.8409996131148
$ws.Range("N21").Value = 40.0997764924509
$ws.Range("O21").Value = 35375.18827407724
$ws.Range("P21").Value = 1936.834205860589
$ws.Range("Q21").Value = 6966.144843695444
$ws.Range("R21").Value = 2695.798690905988
$ws.Range("S21").Value = 239.065753471853
$ws.Range("T21").Value = 724.1259617438004
$ws.Range("U21").Value = 36.7165519629629
$ws.Range("V21").Value = 3419.924652649789
$ws.Range("W21").Value = 275.7823054348158
$ws.Range("X21").Value = 380.8147797481968
$ws.Range("Y21").Value = 20.33850763009366
$ws.Range("Z21").Value = 273.5815285978683
$ws.Range("AA21").Value = 49.8018655397398

